# Apply updated crypto price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that sometimes looks like a plain number
# (e.g. "13.30", "0.0398"). Force those cells to Text format first so
# Excel keeps the exact string (trailing zeros, decimal grouping, etc.)
# instead of silently re-parsing it as a numeric value.

$ws.Range("D2").Value = '64.020.03'
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").Value = '3.151.74'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.85'
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.56'
$ws.Range("E6").Value = '  -2.88%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.145.92'
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("E10").Value = '  -1.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.38'
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("E12").Value = '  -1.53%  '
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("E14").Value = '  -2.23%  '
$ws.Range("D15").Value = '3.668.30'
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("E16").Value = '  +2.46%  '
$ws.Range("D17").Value = '64.107.06'
$ws.Range("E17").Value = '  -1.04%  '
$ws.Range("D18").Value = '3.147.80'
$ws.Range("E18").Value = '  -0.51%  '
$ws.Range("E19").Value = '  -1.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '489.69'
$ws.Range("E20").Value = '  +1.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.68'
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.711'
$ws.Range("E22").Value = '  -1.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.65'
$ws.Range("E23").Value = '  -4.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.75'
$ws.Range("E24").Value = '  +4.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.30'
$ws.Range("E25").Value = '  -3.36%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  -2.60%  '
$ws.Range("E28").Value = '  -4.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.02'
$ws.Range("E29").Value = '  +0.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.06'
$ws.Range("E30").Value = '  -1.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.53'
$ws.Range("E31").Value = '  +3.71%  '
$ws.Range("E32").Value = '  -6.17%  '
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("E34").Value = '  -2.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.10'
$ws.Range("E35").Value = '  -3.10%  '
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.74'
$ws.Range("E37").Value = '  -0.62%  '
$ws.Range("D38").Value = '0.0₃0746'
$ws.Range("E38").Value = '  -4.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.98'
$ws.Range("E39").Value = '  -8.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '436.83'
$ws.Range("E40").Value = '  -5.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0398'
$ws.Range("E41").Value = '  -1.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.119'
$ws.Range("E42").Value = '  -1.47%  '
$ws.Range("E43").Value = '  -1.39%  '
$ws.Range("D44").Value = '2.929.82'
$ws.Range("E44").Value = '  +2.89%  '
$ws.Range("E45").Value = '  -3.97%  '
$ws.Range("E46").Value = '  -6.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.40'
$ws.Range("E47").Value = '  -2.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.95'
$ws.Range("E49").Value = '  -2.64%  '
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.63'
$ws.Range("E51").Value = '  +0.19%  '
